# Normalize the "Recorded By" column (G) entries by rotating the
# comma-separated list of recorder identities left by one position
# wherever it currently starts with a non-"System" identity (i.e.
# dnasr281@gmail.com, ... or system, System, ...), matching the
# canonical ordering used upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known exact remappings of the "Recorded By" text (old -> new),
# derived from rotating the comma separated parts left by one.
$map = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
